$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (e.g. "61.927.72" uses a
# thousands-separator style that is not a valid numeric literal), so force
# the cell format to Text before assigning to avoid Excel auto-converting
# the string into a number.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.995.36'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.427.51'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '410.58'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.25'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.738'
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.140'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '43.73'
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000227'
$ws.Range('E12').Value = '  +17.45%  '
$ws.Range('E13').Value = '  +5.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.966.99'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.32'
$ws.Range('E16').Value = '  +4.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.429.43'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.40'
$ws.Range('E18').Value = '  +9.09%  '
$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.08'
$ws.Range('E19').Value = '  +3.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '61.944.65'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '509.87'
$ws.Range('E21').Value = '  +31.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '92.53'
$ws.Range('E22').Value = '  +4.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.33'
$ws.Range('E23').Value = '  +4.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.45'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('E25').Value = '  +3.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '34.85'
$ws.Range('E26').Value = '  +8.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.25'
$ws.Range('E27').Value = '  +9.04%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.78'
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.67'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.19'
$ws.Range('E30').Value = '  +3.39%  '
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '42.08'
$ws.Range('E33').Value = '  -4.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '59.73'
$ws.Range('E34').Value = '  +14.01%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.997'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  +5.25%  '
$ws.Range('E39').Value = '  +3.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.75'
$ws.Range('E40').Value = '  +18.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '147.65'
$ws.Range('E41').Value = '  +4.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.12'
$ws.Range('E42').Value = '  +7.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.94'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.33'
$ws.Range('E45').Value = '  +8.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.83'
$ws.Range('E46').Value = '  +0.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.35'
$ws.Range('E47').Value = '  +21.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.24'
$ws.Range('E48').Value = '  +5.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '120.91'
$ws.Range('E49').Value = '  +27.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.146'
$ws.Range('E50').Value = '  +18.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.140.84'
$ws.Range('E51').Value = '  +1.01%  '
